$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 00:59"

# Update per-country rows: two country pairs were re-ordered (Barein/Venezuela,
# Montenegro/Tunez, Jamaica/Eslovenia), and case totals were refreshed across
# the table for this snapshot.
$ws.Cells.Item(4, 2).Value = 6706522
$ws.Cells.Item(4, 3).Value = 29921
$ws.Cells.Item(4, 4).Value = 3968885
$ws.Cells.Item(4, 5).Value = 2539163
$ws.Cells.Item(4, 7).Value = 346
$ws.Cells.Item(4, 8).Value = 198474

$ws.Cells.Item(6, 2).Value = 4330455
$ws.Cells.Item(6, 3).Value = 14597
$ws.Cells.Item(6, 4).Value = 3573958
$ws.Cells.Item(6, 5).Value = 624872
$ws.Cells.Item(6, 7).Value = 351
$ws.Cells.Item(6, 8).Value = 131625

$ws.Cells.Item(9, 2).Value = 716319
$ws.Cells.Item(9, 3).Value = 7355
$ws.Cells.Item(9, 4).Value = 599385
$ws.Cells.Item(9, 5).Value = 94010
$ws.Cells.Item(9, 7).Value = 190
$ws.Cells.Item(9, 8).Value = 22924

$ws.Cells.Item(24, 2).Value = 261298
$ws.Cells.Item(24, 3).Value = 752
$ws.Cells.Item(24, 4).Value = 235700
$ws.Cells.Item(24, 5).Value = 16170
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 9428

$ws.Cells.Item(29, 2).Value = 136659
$ws.Cells.Item(29, 3).Value = 518
$ws.Cells.Item(29, 4).Value = 120430
$ws.Cells.Item(29, 5).Value = 7058
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 9171

$ws.Cells.Item(37, 2).Value = 101009
$ws.Cells.Item(37, 3).Value = 153
$ws.Cells.Item(37, 4).Value = 84161
$ws.Cells.Item(37, 5).Value = 11200
$ws.Cells.Item(37, 7).Value = 21
$ws.Cells.Item(37, 8).Value = 5648

$ws.Cells.Item(45, 2).Value = 81909
$ws.Cells.Item(45, 3).Value = 251
$ws.Cells.Item(45, 4).Value = 70927
$ws.Cells.Item(45, 5).Value = 8025
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 2957

$ws.Cells.Item(47, 2).Value = 75218
$ws.Cells.Item(47, 3).Value = 674
$ws.Cells.Item(47, 4).Value = 66899
$ws.Cells.Item(47, 5).Value = 6880
$ws.Cells.Item(47, 7).Value = 16
$ws.Cells.Item(47, 8).Value = 1439

$ws.Cells.Item(53, 1).Value = "Barein"
$ws.Cells.Item(53, 2).Value = 60307
$ws.Cells.Item(53, 3).Value = 721
$ws.Cells.Item(53, 4).Value = 53681
$ws.Cells.Item(53, 5).Value = 6414
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 212

$ws.Cells.Item(54, 1).Value = "Venezuela"
$ws.Cells.Item(54, 2).Value = 59630
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 47729
$ws.Cells.Item(54, 5).Value = 11424
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 477

$ws.Cells.Item(56, 2).Value = 56256
$ws.Cells.Item(56, 3).Value = 79
$ws.Cells.Item(56, 4).Value = 44152
$ws.Cells.Item(56, 5).Value = 11022
$ws.Cells.Item(56, 7).Value = 4
$ws.Cells.Item(56, 8).Value = 1082

$ws.Cells.Item(64, 2).Value = 44881
$ws.Cells.Item(64, 3).Value = 53
$ws.Cells.Item(64, 4).Value = 40922
$ws.Cells.Item(64, 5).Value = 2896
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 1063

$ws.Cells.Item(68, 2).Value = 36157
$ws.Cells.Item(68, 3).Value = 188
$ws.Cells.Item(68, 4).Value = 23067
$ws.Cells.Item(68, 5).Value = 12468
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 622

$ws.Cells.Item(81, 2).Value = 20167
$ws.Cells.Item(81, 3).Value = 158
$ws.Cells.Item(81, 4).Value = 18837
$ws.Cells.Item(81, 5).Value = 915
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 415

$ws.Cells.Item(84, 2).Value = 17918
$ws.Cells.Item(84, 3).Value = 27
$ws.Cells.Item(84, 4).Value = 12767
$ws.Cells.Item(84, 5).Value = 4431
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 720

$ws.Cells.Item(93, 2).Value = 12154
$ws.Cells.Item(93, 3).Value = 75
$ws.Cells.Item(93, 4).Value = 10371
$ws.Cells.Item(93, 5).Value = 1518
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 265

$ws.Cells.Item(108, 1).Value = "Montenegro"
$ws.Cells.Item(108, 2).Value = 6712
$ws.Cells.Item(108, 3).Value = 182
$ws.Cells.Item(108, 4).Value = 4507
$ws.Cells.Item(108, 5).Value = 2085
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 120

$ws.Cells.Item(109, 1).Value = "Tunez"
$ws.Cells.Item(109, 2).Value = 6635
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 1991
$ws.Cells.Item(109, 5).Value = 4537
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 107

$ws.Cells.Item(110, 2).Value = 5690
$ws.Cells.Item(110, 3).Value = 12
$ws.Cells.Item(110, 4).Value = 3731
$ws.Cells.Item(110, 5).Value = 1782
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 177

$ws.Cells.Item(125, 1).Value = "Jamaica"
$ws.Cells.Item(125, 2).Value = 3771
$ws.Cells.Item(125, 3).Value = 148
$ws.Cells.Item(125, 4).Value = 1149
$ws.Cells.Item(125, 5).Value = 2580
$ws.Cells.Item(125, 7).Value = 2
$ws.Cells.Item(125, 8).Value = 42

$ws.Cells.Item(126, 1).Value = "Eslovenia"
$ws.Cells.Item(126, 2).Value = 3702
$ws.Cells.Item(126, 3).Value = 99
$ws.Cells.Item(126, 4).Value = 2730
$ws.Cells.Item(126, 5).Value = 837
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 135

$ws.Cells.Item(138, 2).Value = 3042
$ws.Cells.Item(138, 3).Value = 49
$ws.Cells.Item(138, 4).Value = 772
$ws.Cells.Item(138, 5).Value = 2217
$ws.Cells.Item(138, 7).Value = 2
$ws.Cells.Item(138, 8).Value = 53

$ws.Cells.Item(153, 2).Value = 1853
$ws.Cells.Item(153, 3).Value = 41
$ws.Cells.Item(153, 4).Value = 1215
$ws.Cells.Item(153, 5).Value = 582
$ws.Cells.Item(153, 7).Value = 2
$ws.Cells.Item(153, 8).Value = 56
